$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Wins, Losses, Ties in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered/top) by
# copying the format from an existing header cell onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record for every data row (2-48): 78 wins, 84 losses, 0 ties.
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 78  # AD
    $ws.Cells.Item($r, 31).Value = 84  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
